$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing B,C,D values for rows 2-34
$ws.Range("B2").Value = 3.549002040985112
$ws.Range("C2").Value = 5.920910316106079
$ws.Range("D2").Value = 0.5046761629203652
$ws.Range("B3").Value = 4.701567203995966
$ws.Range("C3").Value = 2.128579211505883
$ws.Range("D3").Value = 0.2347006499076568
$ws.Range("B4").Value = 24.89704965903647
$ws.Range("C4").Value = 3.591925060815581
$ws.Range("D4").Value = 0.2393622870667433
$ws.Range("B5").Value = 25.40763435287513
$ws.Range("C5").Value = 5.763896027392959
$ws.Range("D5").Value = 0.473410923521888
$ws.Range("B6").Value = 27.7155993800383
$ws.Range("C6").Value = 3.924989885302224
$ws.Range("D6").Value = 0.4495565034980107
$ws.Range("B7").Value = 28.70589739754918
$ws.Range("C7").Value = 6.547303912417824
$ws.Range("D7").Value = 0.363691940966677
$ws.Range("B8").Value = 28.9474090844941
$ws.Range("C8").Value = 3.594225935587229
$ws.Range("D8").Value = 0.5614032296489461
$ws.Range("B9").Value = 29.28060521935764
$ws.Range("C9").Value = 5.941376481640827
$ws.Range("D9").Value = 0.5599093400145091
$ws.Range("B10").Value = 29.79888294429146
$ws.Range("C10").Value = 8.570841295250922
$ws.Range("D10").Value = 0.5864776022267965
$ws.Range("B11").Value = 31.99535133286131
$ws.Range("C11").Value = 2.685758212745861
$ws.Range("D11").Value = 0.4504252297610689
$ws.Range("B12").Value = 32.28306324382945
$ws.Range("C12").Value = 6.40341971296473
$ws.Range("D12").Value = 0.3608168959997716
$ws.Range("B13").Value = 32.80251407884936
$ws.Range("C13").Value = 4.652788486339528
$ws.Range("D13").Value = 0.5364662361080966
$ws.Range("B14").Value = 34.17173375263295
$ws.Range("C14").Value = 2.656354586914867
$ws.Range("D14").Value = 0.2461498720556617
$ws.Range("B15").Value = 38.65534588369502
$ws.Range("C15").Value = 3.376219820481659
$ws.Range("D15").Value = 0.2846923715809954
$ws.Range("B16").Value = 39.3969847908689
$ws.Range("C16").Value = 5.847797333772785
$ws.Range("D16").Value = 0.610551657851126
$ws.Range("B17").Value = 40.15245631420453
$ws.Range("C17").Value = 5.003571061324997
$ws.Range("D17").Value = 0.4275109575782723
$ws.Range("B18").Value = 40.39704074563207
$ws.Range("C18").Value = 3.011084607311341
$ws.Range("D18").Value = 0.5602818326300637
$ws.Range("B19").Value = 44.68059995034834
$ws.Range("C19").Value = 8.267586281539778
$ws.Range("D19").Value = 0.7541075925556845
$ws.Range("B20").Value = 44.8731684069034
$ws.Range("C20").Value = 7.744353823533627
$ws.Range("D20").Value = 0.5838520137404215
$ws.Range("B21").Value = 49.53836316326489
$ws.Range("C21").Value = 5.037443392572202
$ws.Range("D21").Value = 0.5192344730762715
$ws.Range("B22").Value = 53.90232659544064
$ws.Range("C22").Value = 2.601411237540201
$ws.Range("D22").Value = 0.4393815164804102
$ws.Range("B23").Value = 54.3267988307084
$ws.Range("C23").Value = 3.726762465324253
$ws.Range("D23").Value = 0.5958474913229082
$ws.Range("B24").Value = 55.53767390222259
$ws.Range("C24").Value = 8.227466583044626
$ws.Range("D24").Value = 0.5498064013963135
$ws.Range("B25").Value = 55.6124258706951
$ws.Range("C25").Value = 1.371158548615634
$ws.Range("D25").Value = 0.3261154478608849
$ws.Range("B26").Value = 56.07269643889632
$ws.Range("C26").Value = 3.128966456147886
$ws.Range("D26").Value = 0.3457104690576793
$ws.Range("B27").Value = 60.40449516779562
$ws.Range("C27").Value = 2.900427098804642
$ws.Range("D27").Value = 0.233726525876494
$ws.Range("B28").Value = 62.77820245342685
$ws.Range("C28").Value = 4.022391567723837
$ws.Range("D28").Value = 0.1935328123797836
$ws.Range("B29").Value = 67.10348550367897
$ws.Range("C29").Value = 2.088993628895523
$ws.Range("D29").Value = 0.2712829443666336
$ws.Range("B30").Value = 68.41838734751978
$ws.Range("C30").Value = 3.180682948667949
$ws.Range("D30").Value = 0.2132067120687068
$ws.Range("B31").Value = 80.21484724568614
$ws.Range("C31").Value = 8.162889871581445
$ws.Range("D31").Value = 0.5587576147131091
$ws.Range("B32").Value = 80.62343184066046
$ws.Range("C32").Value = 6.072643420591841
$ws.Range("D32").Value = 0.6281436017132139
$ws.Range("B33").Value = 81.12651233713237
$ws.Range("C33").Value = 6.058105028906368
$ws.Range("D33").Value = 0.487615366008791
$ws.Range("B34").Value = 84.56460148733915
$ws.Range("C34").Value = 3.200359441969615
$ws.Range("D34").Value = 0.4500313884561341

# Copy format (style) from A34 down to A35:A46 for the new rows
$ws.Range("A34").Copy()
$ws.Range("A35:A46").PasteSpecial(-4122)

# Fill in new rows 35-46 (A=33..44, with B,C,D values)
$ws.Range("A35").Value = 33
$ws.Range("B35").Value = 86.26027492089497
$ws.Range("C35").Value = 2.905659684310158
$ws.Range("D35").Value = 0.5112723520372343
$ws.Range("A36").Value = 34
$ws.Range("B36").Value = 88.11297313672991
$ws.Range("C36").Value = 2.186518346121565
$ws.Range("D36").Value = 0.3934583626101703
$ws.Range("A37").Value = 35
$ws.Range("B37").Value = 88.40250933103228
$ws.Range("C37").Value = 5.670766916117866
$ws.Range("D37").Value = 0.6545677563386395
$ws.Range("A38").Value = 36
$ws.Range("B38").Value = 89.13628026817405
$ws.Range("C38").Value = 10.07173058047172
$ws.Range("D38").Value = 0.388034103967287
$ws.Range("A39").Value = 37
$ws.Range("B39").Value = 90.37011428313143
$ws.Range("C39").Value = 6.701078650276814
$ws.Range("D39").Value = 0.6600837257312172
$ws.Range("A40").Value = 38
$ws.Range("B40").Value = 92.18553353998912
$ws.Range("C40").Value = 4.885797870027363
$ws.Range("D40").Value = 0.2321479262673757
$ws.Range("A41").Value = 39
$ws.Range("B41").Value = 95.21343659189858
$ws.Range("C41").Value = 4.426641546842181
$ws.Range("D41").Value = 0.8961203784555977
$ws.Range("A42").Value = 40
$ws.Range("B42").Value = 95.80760528224833
$ws.Range("C42").Value = 7.430252468662544
$ws.Range("D42").Value = 0.8438951043189388
$ws.Range("A43").Value = 41
$ws.Range("B43").Value = 97.05764603712879
$ws.Range("C43").Value = 6.884909731562365
$ws.Range("D43").Value = 0.5527859675763734
$ws.Range("A44").Value = 42
$ws.Range("B44").Value = 97.29157361717638
$ws.Range("C44").Value = 4.071021596953989
$ws.Range("D44").Value = 0.5889553370169347
$ws.Range("A45").Value = 43
$ws.Range("B45").Value = 97.61217982379853
$ws.Range("C45").Value = 5.968410396222864
$ws.Range("D45").Value = 0.6597982972340297
$ws.Range("A46").Value = 44
$ws.Range("B46").Value = 97.81717229296741
$ws.Range("C46").Value = 4.262582811051101
$ws.Range("D46").Value = 0.2204737389845027
